# Add the header row (Type, Number, Name, Notes) to the empty sheet,
# widen the Notes column, and leave column D selected — matching the
# state the workbook was in when it should have been saved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Type"
$ws.Range("B1").Value = "Number"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Notes"

$ws.Columns.Item(4).ColumnWidth = 112

$ws.Columns.Item(4).Select()
